# Applies the "Updated cryptos list" data refresh to the cryptocurrency
# tracking worksheet. For each changed cell we briefly force a "Text"
# number format before assigning the value so that numeric-looking
# strings (e.g. "1.00", "40.391.45") are preserved verbatim as text
# instead of being normalized into numbers, then we restore the cell's
# original (default/"Normal") style so no extra formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "40.391.45"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -3.23%  "
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.358.79"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -4.64%  "
$cell.Style = "Normal"

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "309.43"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -2.89%  "
$cell.Style = "Normal"

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "85.93"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -7.65%  "
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.524"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -5.27%  "
$cell.Style = "Normal"

# Row 8
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell.Style = "Normal"

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.488"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -5.58%  "
$cell.Style = "Normal"

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0833"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -5.05%  "
$cell.Style = "Normal"

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "30.25"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -8.51%  "
$cell.Style = "Normal"

# Row 12
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -1.52%  "
$cell.Style = "Normal"

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "2.734.46"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -4.30%  "
$cell.Style = "Normal"

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.48"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -6.28%  "
$cell.Style = "Normal"

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "14.83"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -4.99%  "
$cell.Style = "Normal"

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.367.08"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -4.62%  "
$cell.Style = "Normal"

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.751"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -6.39%  "
$cell.Style = "Normal"

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "40.381.70"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -3.16%  "
$cell.Style = "Normal"

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0900"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -5.30%  "
$cell.Style = "Normal"

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "6.06"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -6.75%  "
$cell.Style = "Normal"

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "68.13"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -4.26%  "
$cell.Style = "Normal"

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.61"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -6.42%  "
$cell.Style = "Normal"

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "231.39"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -4.26%  "
$cell.Style = "Normal"

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.61"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -5.38%  "
$cell.Style = "Normal"

# Row 25
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  +0.10%  "
$cell.Style = "Normal"

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.79"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -8.52%  "
$cell.Style = "Normal"

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "23.37"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -7.04%  "
$cell.Style = "Normal"

# Row 28
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -3.08%  "
$cell.Style = "Normal"

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.22"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -5.47%  "
$cell.Style = "Normal"

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "33.37"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -9.65%  "
$cell.Style = "Normal"

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "151.63"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -4.16%  "
$cell.Style = "Normal"

# Row 32
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +0.11%  "
$cell.Style = "Normal"

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.14"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -6.90%  "
$cell.Style = "Normal"

# Row 34
$cell = $ws.Range("B34")
$cell.NumberFormat = "@"
$cell.Value = "WEMIXToken"
$cell.Style = "Normal"
$cell = $ws.Range("C34")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.44"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -4.77%  "
$cell.Style = "Normal"

# Row 35
$cell = $ws.Range("B35")
$cell.NumberFormat = "@"
$cell.Value = "Hedera"
$cell.Style = "Normal"
$cell = $ws.Range("C35")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.0722"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -5.48%  "
$cell.Style = "Normal"

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.113"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -2.91%  "
$cell.Style = "Normal"

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.72"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -6.69%  "
$cell.Style = "Normal"

# Row 38
$cell = $ws.Range("B38")
$cell.NumberFormat = "@"
$cell.Value = "Celestia"
$cell.Style = "Normal"
$cell = $ws.Range("C38")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "15.60"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -10.14%  "
$cell.Style = "Normal"

# Row 39
$cell = $ws.Range("B39")
$cell.NumberFormat = "@"
$cell.Value = "Kaspa"
$cell.Style = "Normal"
$cell = $ws.Range("C39")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0980"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -6.33%  "
$cell.Style = "Normal"

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.68"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -10.05%  "
$cell.Style = "Normal"

# Row 41
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -4.70%  "
$cell.Style = "Normal"

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.79"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -5.79%  "
$cell.Style = "Normal"

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.941.80"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -2.74%  "
$cell.Style = "Normal"

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0266"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -6.41%  "
$cell.Style = "Normal"

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "17.43"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -9.01%  "
$cell.Style = "Normal"

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "9.39"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -0.93%  "
$cell.Style = "Normal"

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.67"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -10.38%  "
$cell.Style = "Normal"

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.597.14"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -4.39%  "
$cell.Style = "Normal"

# Row 49
$cell = $ws.Range("B49")
$cell.NumberFormat = "@"
$cell.Value = "Aave"
$cell.Style = "Normal"
$cell = $ws.Range("C49")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "92.28"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -6.11%  "
$cell.Style = "Normal"

# Row 50
$cell = $ws.Range("B50")
$cell.NumberFormat = "@"
$cell.Value = "BitcoinSV"
$cell.Style = "Normal"
$cell = $ws.Range("C50")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "71.81"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -6.10%  "
$cell.Style = "Normal"

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "49.85"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -5.21%  "
$cell.Style = "Normal"

